# ============================================================================
# Edit: add 2022-Q3 data
#  1. Insert a new "2022-Q3" worksheet right after "总计", populate it with
#     the quarterly fund-holding breakdown, and push the existing quarter
#     sheets (2022-Q2 .. 2021-Q1) one position to the right.
#  2. Insert a new summary row in "总计" for 2022-Q3 (38 holdings, 11.41
#     亿元), shifting the previously-existing summary rows down by one and
#     renumbering the index column.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# Step 1: "总计" (summary) sheet - insert the new 2022-Q3 row at the top of
# the data (row 2), shifting everything else down.
# ----------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()
# The freshly-inserted row can inherit stray formatting from the Insert
# operation; clear it so the data cells come back to the default (unstyled)
# look used by every other data row in this sheet.
$summary.Range("B2:D2").ClearFormats()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 38
$summary.Range("D2").Value = 11.41

# Renumber the index column for the rows that got pushed down (previously
# 0..5 on rows 2..7, now 1..6 on rows 3..8).
for ($i = 3; $i -le 8; $i++) {
    $summary.Range("A$i").Value = ($i - 2)
}

# Make sure the new A2 cell carries the same style as the rest of column A.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

# Data for new sheet "2022-Q3"
$q3Data = @(
    ,('004698','博时军工主题股票','33.59','93.18','6.77','2.2740',6)
    ,('050009','博时新兴成长混合','24.16','92.80','8.75','2.1140',1)
    ,('516150','嘉实中证稀土产业ETF','20.00','99.35','4.18','0.8360',8)
    ,('011162','博时港股通领先趋势混合A','15.21','90.30','5.40','0.8213',7)
    ,('011486','博时创新精选混合A','9.52','93.58','6.37','0.6064',1)
    ,('011756','博时产业优选灵活配置混合A','24.41','72.85','2.40','0.5858',10)
    ,('050010','博时特许价值混合','6.01','92.93','8.43','0.5066',1)
    ,('257020','国联安精选混合','10.39','92.50','4.60','0.4779',3)
    ,('516780','华泰柏瑞中证稀土产业ETF','8.05','98.92','4.10','0.3300',8)
    ,('011592','博时军工主题股票C','4.52','93.18','6.77','0.3060',6)
    ,('005358','东方阿尔法精选灵活配置混合A','3.47','93.91','8.63','0.2995',4)
    ,('050022','博时回报混合','4.17','76.75','5.51','0.2298',2)
    ,('001047','光大保德信国企改革主题股票','2.88','90.88','7.92','0.2281',5)
    ,('011163','博时港股通领先趋势混合C','3.82','90.30','5.40','0.2063',7)
    ,('012082','博时数字经济18个月封闭混合A','5.19','96.61','3.93','0.2040',5)
    ,('006864','国联安核心资产策略混合','4.64','91.24','3.60','0.1670',8)
    ,('014036','博时成长回报混合A','3.81','92.27','4.30','0.1638',3)
    ,('013836','博时时代消费混合A','4.62','93.57','3.46','0.1599',8)
    ,('014325','国联安核心趋势一年持有混合A','3.74','86.69','3.48','0.1302',9)
    ,('001215','博时沪港深优质企业混合A','2.22','91.96','5.59','0.1241',4)
    ,('159715','易方达中证稀土产业ETF','2.52','98.35','4.14','0.1043',8)
    ,('001463','光大保德信一带一路战略主题混合','1.57','87.51','6.12','0.0961',2)
    ,('159713','富国中证稀土产业ETF','2.25','98.40','4.15','0.0934',8)
    ,('014600','博时回报严选混合A','0.92','92.53','8.77','0.0807',2)
    ,('011487','博时创新精选混合C','1.03','93.58','6.37','0.0656',1)
    ,('005359','东方阿尔法精选灵活配置混合C','0.58','93.91','8.63','0.0501',4)
    ,('014037','博时成长回报混合C','0.80','92.27','4.30','0.0344',3)
    ,('012696','同泰数字经济主题股票A','0.99','94.09','2.81','0.0278',5)
    ,('009317','金信核心竞争力灵活配置混合','0.22','89.48','8.94','0.0197',3)
    ,('011757','博时产业优选灵活配置混合C','0.73','72.85','2.40','0.0175',10)
    ,('014326','国联安核心趋势一年持有混合C','0.38','86.69','3.48','0.0132',9)
    ,('012083','博时数字经济18个月封闭混合C','0.32','96.61','3.93','0.0126',5)
    ,('012697','同泰数字经济主题股票C','0.38','94.09','2.81','0.0107',5)
    ,('002555','博时沪港深优质企业混合C','0.08','91.96','5.59','0.0045',4)
    ,('014601','博时回报严选混合C','0.04','92.53','8.77','0.0035',2)
    ,('013837','博时时代消费混合C','0.10','93.57','3.46','0.0035',8)
    ,('004402','金信民旺债券C','0.09','23.60','1.63','0.0015',5)
    ,('004222','金信民旺债券A','0.08','23.60','1.63','0.0013',5)
)

# ----------------------------------------------------------------------------
# Step 2: create the new "2022-Q3" worksheet right after "总计" and fill in
# the fund-holding table (same layout/style as the other quarter sheets).
# ----------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $summary)
$newSheet.Name = "2022-Q3"

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$row = 2
foreach ($rec in $q3Data) {
    $newSheet.Range("A$row").Value = ($row - 2)

    $newSheet.Range("B$row").NumberFormat = "@"
    $newSheet.Range("B$row").Value = $rec[0]

    $newSheet.Range("C$row").Value = $rec[1]

    $newSheet.Range("D$row").NumberFormat = "@"
    $newSheet.Range("D$row").Value = $rec[2]

    $newSheet.Range("E$row").NumberFormat = "@"
    $newSheet.Range("E$row").Value = $rec[3]

    $newSheet.Range("F$row").NumberFormat = "@"
    $newSheet.Range("F$row").Value = $rec[4]

    $newSheet.Range("G$row").NumberFormat = "@"
    $newSheet.Range("G$row").Value = $rec[5]

    $newSheet.Range("H$row").Value = $rec[6]

    $row++
}
$lastRow = $row - 1

# Match the header/index-column style used by every other quarter sheet by
# copying the formatting straight from the sheet that used to be "2022-Q2"
# (now pushed one slot to the right by the insert above) rather than
# reconstructing it through the Font/Border object model (which would add a
# slightly different - if visually equivalent - style entry).
$styleSource = $wb.Worksheets.Item(3)
$styleSource.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats
$styleSource.Range("A2").Copy()
$newSheet.Range("A2:A$lastRow").PasteSpecial(-4122)   # xlPasteFormats

$newSheet.Activate()
